# Update the "Förändrad" date column (C) for all data rows (2-196)
# from serial date 45180 (2023-09-11) to 45181 (2023-09-12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 196
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45180) {
        $cell.Value2 = 45181
    }
}
